$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1133
    3  = 842
    6  = 1111
    9  = 7712
    12 = 379
    13 = 152
    14 = 418
    16 = 7900
    18 = 1374
    30 = 1149
    35 = 45
    36 = 79
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
